$d = $word.ActiveDocument
$d.TrackRevisions = $false

# 1. Update the timestamp paragraph near the top of the document.
$d.Content.Find.Execute("6/1/2023 12:39:40 AM", $true, $false, $false, $false, $false, `
    $true, 1, $false, "6/1/2023 3:49:24 AM", 2) | Out-Null

Write-Output "Step1 done"

# Paragraph indices (1-based, stable at this point in the script):
#   53: "PREVENTION SECURITY SYSTEM: ANY TREASON TOES;"
#   54: "PREVENTION SECURITY SYSTEM: ANY TWIDDLE TOES;"
#   55: "PREVENTION SECURITY SYSTEM: ANY TWINKLE TOES;"
#   56: "PREVENTION SECURITY SYSTEM: ANY UNETIQUETTE;"   (carries w:lastRenderedPageBreak)

# 2. Insert a brand-new paragraph (a duplicate of the still-intact "TWINKLE"
#    paragraph) just before the "UNETIQUETTE" paragraph, and move the
#    w:lastRenderedPageBreak marker onto it.
$pUnetiquette = $d.Paragraphs(56)
$insertPoint = $pUnetiquette.Range.Duplicate()
$insertPoint.Collapse(1)
$insertPoint.InsertBefore("X`r")

Write-Output "After InsertBefore, paragraph count: $($d.Paragraphs.Count)"
Write-Output "p56: [$($d.Paragraphs(56).Range.Text)]"
Write-Output "p57: [$($d.Paragraphs(57).Range.Text)]"

$newParaRange = $d.Paragraphs(56).Range.Duplicate()
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="720"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>PREVENTION SECURITY SYSTEM</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">ANY </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">TWINKLE </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>TOES</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>;</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$newParaRange.InsertXML($newParaXml)

Write-Output "After InsertXML, paragraph count: $($d.Paragraphs.Count)"
for ($i = 53; $i -le 58; $i++) {
    Write-Output "$i : [$($d.Paragraphs($i).Range.Text)]"
}

# 3. Strip the w:lastRenderedPageBreak from the (now shifted) "UNETIQUETTE"
#    paragraph -- replace its contents with an equivalent copy lacking the
#    page-break marker.
$pUnetiquette2 = $d.Paragraphs(57)
Write-Output "p57 check: [$($pUnetiquette2.Range.Text)]"
$unetRange = $pUnetiquette2.Range.Duplicate()
$unetXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="720"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>PREVENTION SECURITY SYSTEM</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">ANY </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>UNETIQUETTE</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>;</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$unetRange.InsertXML($unetXml)

Write-Output "After UNETIQUETTE fix, paragraph count: $($d.Paragraphs.Count)"
for ($i = 53; $i -le 58; $i++) {
    Write-Output "$i : [$($d.Paragraphs($i).Range.Text)]"
}

# 4. "TWINKLE " -> "TWIDDLE " in paragraph 55 (scoped Find avoids touching
#    the freshly-inserted duplicate paragraph 56).
$p55 = $d.Paragraphs(55)
$p55.Range.Find.Execute("TWINKLE ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "TWIDDLE ", 2) | Out-Null
Write-Output "p55 after: [$($d.Paragraphs(55).Range.Text)]"

# 5. "TWIDDLE " -> "TREASON " in paragraph 54.
$p54 = $d.Paragraphs(54)
$p54.Range.Find.Execute("TWIDDLE ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "TREASON ", 2) | Out-Null
Write-Output "p54 after: [$($d.Paragraphs(54).Range.Text)]"

for ($i = 53; $i -le 58; $i++) {
    Write-Output "$i : [$($d.Paragraphs($i).Range.Text)]"
}
